# Auto-generated edit script: refresh cryptos list values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.123.67"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "1.835.53"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "'244.05"
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("D6").Value = "'0.6285"
$ws.Range("E6").Value = "  +0.93%  "
$ws.Range("D7").Value = "'1.005"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").Value = "'0.07487"
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("D9").Value = "'0.2924"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("D10").Value = "'23.08"
$ws.Range("E10").Value = "  +2.17%  "
$ws.Range("D11").Value = "'0.07734"
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "1.829.37"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").Value = "'4.985"
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").Value = "'0.6690"
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("D15").Value = "'82.45"
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").Value = "'0.000009345"
$ws.Range("E16").Value = "  -6.22%  "
$ws.Range("D17").Value = "'5.993"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").Value = "29.142.71"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").Value = "2.077.66"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("E20").Value = "  +2.02%  "
$ws.Range("D21").Value = "'223.08"
$ws.Range("E21").Value = "  -1.13%  "
$ws.Range("D22").Value = "'1.006"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("D23").Value = "'7.126"
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").Value = "'1.005"
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'160.35"
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").Value = "'0.1396"
$ws.Range("E26").Value = "  +1.49%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'8.502"
$ws.Range("E27").Value = "  +0.71%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'17.95"
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'1.502"
$ws.Range("E29").Value = "  +0.94%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.05895"
$ws.Range("E30").Value = "  +13.34%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.155"
$ws.Range("E31").Value = "  +2.12%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'4.063"
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").Value = "'1.206"
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.7493"
$ws.Range("E34").Value = "  +1.48%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'1.848"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.138"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.683"
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.228.66"
$ws.Range("E38").Value = "  -1.46%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.767"
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01792"
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'6.558"
$ws.Range("E41").Value = "  +2.93%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.8946"
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "'1.007"
$ws.Range("E43").Value = "  +0.55%  "
$ws.Range("B44").Value = "XinFinNetwork"
$ws.Range("C44").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D44").Value = "'0.08139"
$ws.Range("E44").Value = "  +20.23%  "
$ws.Range("D45").Value = "'102.15"
$ws.Range("E45").Value = "  +0.62%  "
$ws.Range("D46").Value = "1.978.58"
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'65.58"
$ws.Range("E47").Value = "  +2.36%  "
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.5106"
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").Value = "'0.4068"
$ws.Range("E50").Value = "  +1.24%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'9.028"
$ws.Range("E51").Value = "  +1.75%  "
